# ECSE 425 Final Report — "Optimizations" section edits
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the empty "BodyText" paragraph that sits between the end of the
#    Testing/Evaluation section and the "Optimizations" Heading1 paragraph,
#    so "Optimizations" immediately follows the preceding paragraph.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("design problems that lead to timing errors. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $anchorPara = $d.Range($r1.End, $r1.End).Paragraphs(1)
    $emptyPara = $anchorPara.Next()
    if ($emptyPara.Range.Text.Trim() -eq "") {
        $delRange = $d.Range($emptyPara.Range.Start, $emptyPara.Range.End)
        $delRange.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2) "There are a few ways to optimize a pipelined processor. For our
#    project, we chose..." -> insert " including caching, branch prediction,
#    and early branch detection" right after "processor" (before the period).
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("pipelined processor. For our project, we ch", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $insertPos = $r2.Start + ("pipelined processor".Length)
    $ins2 = $d.Range($insertPos, $insertPos)
    $ins2.InsertAfter(" including caching, branch prediction, and early branch detection")
}

# ---------------------------------------------------------------------------
# 3) "This means that branch instructions are detected..." -> "branching
#    instructions...". Append a new explanatory sentence, and relocate the
#    document's lone "_GoBack" bookmark to sit right after that new sentence
#    (this is where the last edit happened).
# ---------------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute("This means that branch instructions are detected and resolved at the ID stage rather than at the EX stage of the pipeline.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $branchEnd = $r3.Start + ("This means that branch".Length)
    $insIng = $d.Range($branchEnd, $branchEnd)
    $insIng.InsertAfter("ing")

    $newSentence = " This is possible because as soon as the decoder has finished interpreting a branch instruction, it has extracted the target address from the machine code. Therefore, instead of waiting until the EX stage to update the program counter, we can do this immediately following instruction decode."
    $sentenceEnd = $r3.End + "ing".Length
    $insNewSentence = $d.Range($sentenceEnd, $sentenceEnd)
    $insNewSentence.InsertAfter($newSentence)

    $bookmarkPos = $sentenceEnd + $newSentence.Length
    $bmRange = $d.Range($bookmarkPos, $bookmarkPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------------
# 4) Split "... with our early branch detection implemented. We then
#    observed..." so the word "implemented" (and the remainder of the
#    paragraph) starts a new run carrying a lastRenderedPageBreak marker —
#    reflecting where the added text above now pushes the page break to.
# ---------------------------------------------------------------------------
$r4 = $d.Content
$tailPhrase = " implemented. We then observed and compared the delay of our program counter updates. After extensive debugging of our branching mechanism we verified that early branch detection speeds up branch predictions by two clock cycles."
$found4 = $r4.Find.Execute($tailPhrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $target4 = $d.Range($r4.Start, $r4.End)
    $xml4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>implemented. We then observed and compared the delay of our program counter updates</w:t></w:r><w:r><w:t>. After extensive debugging of our branching mechanism</w:t></w:r><w:r><w:t xml:space="preserve"> we</w:t></w:r><w:r><w:t xml:space="preserve"> verified that early branch detection speeds up branch predictions by two clock cycles.</w:t></w:r></w:p>'
    $target4.InsertXML($xml4)
}

# ---------------------------------------------------------------------------
# 5) Remove the now-stale lastRenderedPageBreak that used to sit before
#    "Speedup" (the page break moved earlier in the document, see step 4).
# ---------------------------------------------------------------------------
$r5 = $d.Content
$found5 = $r5.Find.Execute("Speedup", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found5) {
    $para5 = $d.Range($r5.Start, $r5.Start).Paragraphs(1)
    $target5 = $d.Range($para5.Range.Start, $para5.Range.End)
    $xml5 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:i/></w:rPr><w:t>Speedup</w:t></w:r><w:r><w:rPr><w:i/><w:vertAlign w:val="subscript"/></w:rPr><w:t>overall</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> = </w:t></w:r><w:r><w:t>1/[(1-</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>F)+(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">F/S)]  </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(1)</w:t></w:r></w:p>'
    $target5.InsertXML($xml5)
}

# ---------------------------------------------------------------------------
# 6) "...naming convention prescribed by your conference... highlight all of
#    the contents..." -> rewrap runs and flag "all of" with gramStart/gramEnd
#    proofing marks (as Word's grammar checker would). The old bookmark
#    location here is subsumed since the bookmark was relocated in step 3.
# ---------------------------------------------------------------------------
$r6 = $d.Content
$found6 = $r6.Find.Execute("After the text edit has been completed", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found6) {
    $para6 = $d.Range($r6.Start, $r6.Start).Paragraphs(1)
    $target6 = $d.Range($para6.Range.Start, $para6.Range.End)
    $xml6 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">After the text edit has been completed, the paper is ready for the template. Duplicate the template file by using the Save As command, and use the naming convention prescribed by </w:t></w:r><w:r><w:t xml:space="preserve">your conference for the name of your paper. In this newly created file, highlight </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">all </w:t></w:r><w:r><w:t>of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the contents and import your </w:t></w:r><w:r><w:t>prepared text file. You are now ready to style your paper; use the scroll down window on the left of the MS Word Formatting toolbar.</w:t></w:r></w:p>'
    $target6.InsertXML($xml6)
}

Write-Output "All steps completed"
